$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8383
$ws1.Range("F3").Value = 7940
$ws1.Range("F4").Value = 131
$ws1.Range("F8").Value = 133
$ws1.Range("F11").Value = 234
$ws1.Range("F12").Value = 717
$ws1.Range("F13").Value = 136
$ws1.Range("F14").Value = 1924
$ws1.Range("F15").Value = 64
$ws1.Range("F20").Value = 13

# Sheet "全部类型" (sheet4) - column F ("想去人数") updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8383
$ws4.Range("F3").Value = 7941
$ws4.Range("F4").Value = 131
$ws4.Range("F8").Value = 133
$ws4.Range("F11").Value = 234
$ws4.Range("F12").Value = 717
$ws4.Range("F13").Value = 136
$ws4.Range("F14").Value = 1924
$ws4.Range("F15").Value = 64
$ws4.Range("F20").Value = 13
